# "7-not" grade sheet: fill in the student's exam scores -> year-end grade
# formula, the student info block, and the class-average formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Year-end grade formula: 1st midterm 20% + 2nd midterm 20% + final 60% ---
# F3 is entered on its own; F4:F9 are filled down together so Excel records
# them as one shared formula (matches the original author's workflow).
$ws.Range("F3").Formula = "=C3*0.2+D3*0.2+E3*0.6"
$ws.Range("F4:F9").Formula = "=C4*0.2+D4*0.2+E4*0.6"

# Row 9 picked up its own one-off cell style (thicker outer border) before the
# formula existed; copy the plain F3:F8 formatting onto it so it matches the
# rest of the column while keeping the formula/result we just wrote.
$ws.Range("F8").Copy()
$ws.Range("F9").PasteSpecial(-4122)   # xlPasteFormats

# --- Student info block (Numara / Ad Soyad / Bölüm) ---
$ws.Range("I4").Value = 20215070019
$ws.Range("I5").Value = "KÜBRA ÇABUK"
$ws.Range("I6").Value = "YBS"

# --- Class average (year-end grade) ---
$ws.Range("C14").Formula = "=AVERAGE(F3:F9)"
$ws.Range("C14").NumberFormat = '_-* #,##0.00_-;\-* #,##0.00_-;_-* "-"??_-;_-@_-'

# --- Leave the selection where the editing session ended ---
$ws.Range("H11").Select()
